# Updated symbol list: refresh Price (D) and Volume(1h) (E) columns for the
# coin rows whose figures changed in this run, preserving the existing
# "text" storage of these cells (they are plain strings, not numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value, taken from the refreshed feed.
$updates = [ordered]@{
    'D2' = '305.77'
    'E2' = '-0.05%'
    'D3' = '35.61'
    'E3' = '-0.59%'
    'D4' = '5.037'
    'E4' = '-0.90%'
    'D5' = '0.08004'
    'E5' = '-0.86%'
    'D6' = '1.919'
    'E6' = '-0.29%'
    'D7' = '7.776'
    'E7' = '0.22%'
    'D8' = '0.9202'
    'D9' = '0.1274'
    'E9' = '-7.19%'
    'D10' = '0.1914'
    'E10' = '0.49%'
    'D11' = '0.09115'
    'E11' = '-0.71%'
    'D12' = '0.03457'
    'E12' = '1.09%'
    'D13' = '0.09840'
    'E13' = '0.04%'
    'E14' = '-0.32%'
    'D15' = '0.006311'
    'E15' = '9.46%'
    'D16' = '3.720'
    'E16' = '2.64%'
    'D17' = '4.160'
    'E17' = '-0.76%'
    'E18' = '12.60%'
    'D19' = '0.3443'
    'E19' = '-0.39%'
    'E20' = '2.82%'
    'D21' = '5.175'
    'E21' = '5.24%'
    'D22' = '0.2604'
    'E22' = '6.52%'
    'D23' = '0.04444'
    'E23' = '0.28%'
    'E24' = '0.97%'
    'D25' = '0.004622'
    'E25' = '-3.84%'
    'D26' = '0.0001252'
    'E26' = '-3.98%'
    'D27' = '0.0004447'
    'E27' = '41.91%'
    'D39' = '0.01947'
    'E39' = '-3.77%'
    'D40' = '0.05364'
    'E40' = '8.89%'
    'D41' = '0.007608'
    'E41' = '-0.41%'
    'E42' = '-1.70%'
    'D43' = '0.1354'
    'E43' = '-1.67%'
    'D44' = '0.002153'
    'E44' = '2.24%'
    'D45' = '0.009900'
    'E45' = '-10.56%'
    'D46' = '0.00006135'
    'E46' = '-4.52%'
    'D47' = '0.00000000751'
    'E47' = '-0.04%'
    'E48' = '0.85%'
    'D49' = '0.001660'
    'E49' = '39.23%'
    'D50' = '0.00002103'
    'E50' = '-0.04%'
    'D51' = '0.0002003'
    'E51' = '-0.04%'
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text storage (no leading quote / no style change) so numeric-
    # looking strings like "305.77" or "-0.05%" stay text, matching the
    # original inline-string cell type instead of being parsed as a number.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
